# Reverse the order of comma-separated "Recorded By" entries in column G.
# Cells that contain a single value (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = 7
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.Contains(",")) {
        $parts = $value -split ","
        $trimmed = @()
        foreach ($part in $parts) {
            $trimmed += $part.Trim()
        }

        $reversed = @()
        for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
            $reversed += $trimmed[$i]
        }

        $newValue = [string]::Join(", ", $reversed)
        $cell.Value = $newValue
    }
}
